$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date values in column B (stored as Excel serial date numbers)
$ws.Range("B2").Value = 41821
$ws.Range("B3").Value = 41835
$ws.Range("B4").Value = 41883
$ws.Range("B5").Value = 41927
$ws.Range("B6").Value = 41958

# Update row heights that Excel recalculated for these rows
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 37.5
$ws.Rows.Item(10).RowHeight = 131.25
$ws.Rows.Item(13).RowHeight = 37.5

# Rows 11 and 15 shrink back to the sheet's default (standard) row height
$ws.Rows.Item(11).AutoFit()
$ws.Rows.Item(15).AutoFit()

# Move the active selection from A1:F6 to C6
$ws.Range("C6").Select()
